$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: swap "Aktivitas 1" -> "Transaksi 2", update dates
$ws.Range("B2").Value = "Transaksi 2"
$ws.Range("C2").Value = 45027
$ws.Range("D2").Value = 45080

# Update row 3: swap "Aktivitas 2" -> "Transaksi 3", update dates
$ws.Range("B3").Value = "Transaksi 3"
$ws.Range("C3").Value = 45261
$ws.Range("D3").Value = 45084

# Remove row 4 (was previously C4/D4 with only styles, no values)
$ws.Range("A4:D4").Delete()

# Update selection to G8
$ws.Range("G8").Select()
